$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column E (weighted_avg_execution_spread_$) -- shifts F:L left to E:K,
# matching the new dimension A1:K28 and dropping that metric column entirely.
$ws.Columns("E").Delete()

# Column H (new) needs to go from width 21 to width 20 to match the refreshed layout
# (21 - 0.83 internal offset = 20.17 -> stored width 21; 20 - 0.83 = 19.17 -> stored width 20)
$ws.Columns("H").ColumnWidth = 19.17

# ---- EURUSD block: refreshed header row 2 and data rows ----
$ws.Cells.Item(2, 2).Value = "count_of_occurrences"
$ws.Cells.Item(2, 3).Value = "percentage_of_occurrences"
$ws.Cells.Item(2, 4).Value = "typical_spread_in_points"
$ws.Cells.Item(2, 5).Value = "volume_weighted_avg_spread_in_USD"
$ws.Cells.Item(2, 6).Value = "PnL_per_lot"
$ws.Cells.Item(2, 7).Value = "total_profit"
$ws.Cells.Item(2, 8).Value = "pct_total_profit"
$ws.Cells.Item(2, 9).Value = "total_volume"
$ws.Cells.Item(2, 10).Value = "pct_total_volume"
$ws.Cells.Item(2, 11).Value = "pct_impact_on_PnL_exec_spread"
$ws.Cells.Item(4, 2).Value = 288
$ws.Cells.Item(4, 3).Value = 36.97047496790758
$ws.Cells.Item(4, 4).Value = 9.183673469387756
$ws.Cells.Item(4, 5).Value = 8.32135451943875
$ws.Cells.Item(4, 6).Value = -6.695116692594337
$ws.Cells.Item(4, 7).Value = -29406158.84
$ws.Cells.Item(4, 8).Value = 27.77707824953054
$ws.Cells.Item(4, 9).Value = 4392180.18
$ws.Cells.Item(4, 10).Value = 37.85304652439557
$ws.Cells.Item(4, 11).Value = 0.1493625945468776
$ws.Cells.Item(5, 2).Value = 176
$ws.Cells.Item(5, 3).Value = 22.59306803594352
$ws.Cells.Item(5, 4).Value = 9.241379310344827
$ws.Cells.Item(5, 5).Value = 8.728223779357057
$ws.Cells.Item(5, 6).Value = -17.23438850857107
$ws.Cells.Item(5, 7).Value = -48573527.59
$ws.Cells.Item(5, 8).Value = 45.88258820420493
$ws.Cells.Item(5, 9).Value = 2818407.37
$ws.Cells.Item(5, 10).Value = 24.28982895262493
$ws.Cells.Item(5, 11).Value = 0.05802352659641372
$ws.Cells.Item(6, 2).Value = 209
$ws.Cells.Item(6, 3).Value = 26.82926829268293
$ws.Cells.Item(6, 4).Value = 10.69230769230769
$ws.Cells.Item(6, 5).Value = 8.144813843788437
$ws.Cells.Item(6, 6).Value = -5.240664324645246
$ws.Cells.Item(6, 7).Value = -15247926.2
$ws.Cells.Item(6, 8).Value = 14.40320177500839
$ws.Cells.Item(6, 9).Value = 2909540.71
$ws.Cells.Item(6, 10).Value = 25.07524175846833
$ws.Cells.Item(6, 11).Value = 0.1908155031600297
$ws.Cells.Item(7, 2).Value = 106
$ws.Cells.Item(7, 3).Value = 13.60718870346598
$ws.Cells.Item(7, 4).Value = 9.142857142857142
$ws.Cells.Item(7, 5).Value = 8.053686994388862
$ws.Cells.Item(7, 6).Value = -8.52074533246008
$ws.Cells.Item(7, 7).Value = -12637225.19
$ws.Cells.Item(7, 8).Value = 11.93713177125614
$ws.Cells.Item(7, 9).Value = 1483112.65
$ws.Cells.Item(7, 10).Value = 12.78188276451118
$ws.Cells.Item(7, 11).Value = 0.1173606252718837

# ---- GBPUSD block: refreshed header row 9 and data rows ----
$ws.Cells.Item(9, 2).Value = "count_of_occurrences"
$ws.Cells.Item(9, 3).Value = "percentage_of_occurrences"
$ws.Cells.Item(9, 4).Value = "typical_spread_in_points"
$ws.Cells.Item(9, 5).Value = "volume_weighted_avg_spread_in_USD"
$ws.Cells.Item(9, 6).Value = "PnL_per_lot"
$ws.Cells.Item(9, 7).Value = "total_profit"
$ws.Cells.Item(9, 8).Value = "pct_total_profit"
$ws.Cells.Item(9, 9).Value = "total_volume"
$ws.Cells.Item(9, 10).Value = "pct_total_volume"
$ws.Cells.Item(9, 11).Value = "pct_impact_on_PnL_exec_spread"
$ws.Cells.Item(11, 2).Value = 297
$ws.Cells.Item(11, 3).Value = 38.12580231065468
$ws.Cells.Item(11, 4).Value = 9.95774647887324
$ws.Cells.Item(11, 5).Value = 11.64557714698581
$ws.Cells.Item(11, 6).Value = -12.10709004729152
$ws.Cells.Item(11, 7).Value = -33600487.26
$ws.Cells.Item(11, 8).Value = 33.29042947458734
$ws.Cells.Item(11, 9).Value = 2775273.59
$ws.Cells.Item(11, 10).Value = 39.95509274347251
$ws.Cells.Item(11, 11).Value = 0.08259623048097428
$ws.Cells.Item(12, 2).Value = 141
$ws.Cells.Item(12, 3).Value = 18.10012836970475
$ws.Cells.Item(12, 4).Value = 12.125
$ws.Cells.Item(12, 5).Value = 12.94790566732876
$ws.Cells.Item(12, 6).Value = -40.74193108309967
$ws.Cells.Item(12, 7).Value = -46853968.36
$ws.Cells.Item(12, 8).Value = 46.42161041366776
$ws.Cells.Item(12, 9).Value = 1150018.35
$ws.Cells.Item(12, 10).Value = 16.55659823828224
$ws.Cells.Item(12, 11).Value = 0.02454473740973005
$ws.Cells.Item(13, 2).Value = 284
$ws.Cells.Item(13, 3).Value = 36.45699614890886
$ws.Cells.Item(13, 4).Value = 12.30769230769231
$ws.Cells.Item(13, 5).Value = 12.12941745316376
$ws.Cells.Item(13, 6).Value = -4.564054699990865
$ws.Cells.Item(13, 7).Value = -11558985.27
$ws.Cells.Item(13, 8).Value = 11.45232153781358
$ws.Cells.Item(13, 9).Value = 2532613.22
$ws.Cells.Item(13, 10).Value = 36.46155696254961
$ws.Cells.Item(13, 11).Value = 0.2191034213507567
$ws.Cells.Item(14, 2).Value = 57
$ws.Cells.Item(14, 3).Value = 7.317073170731707
$ws.Cells.Item(14, 4).Value = 10.33333333333333
$ws.Cells.Item(14, 5).Value = 12.25270450645865
$ws.Cells.Item(14, 6).Value = -18.27156851540661
$ws.Cells.Item(14, 7).Value = -8917931.250000002
$ws.Cells.Item(14, 8).Value = 8.835638573931313
$ws.Cells.Item(14, 9).Value = 488076.94
$ws.Cells.Item(14, 10).Value = 7.026752055695622
$ws.Cells.Item(14, 11).Value = 0.05472983882893243

# ---- USDJPY block: refreshed header row 16 and data rows ----
$ws.Cells.Item(16, 2).Value = "count_of_occurrences"
$ws.Cells.Item(16, 3).Value = "percentage_of_occurrences"
$ws.Cells.Item(16, 4).Value = "typical_spread_in_points"
$ws.Cells.Item(16, 5).Value = "volume_weighted_avg_spread_in_USD"
$ws.Cells.Item(16, 6).Value = "PnL_per_lot"
$ws.Cells.Item(16, 7).Value = "total_profit"
$ws.Cells.Item(16, 8).Value = "pct_total_profit"
$ws.Cells.Item(16, 9).Value = "total_volume"
$ws.Cells.Item(16, 10).Value = "pct_total_volume"
$ws.Cells.Item(16, 11).Value = "pct_impact_on_PnL_exec_spread"
$ws.Cells.Item(18, 2).Value = 311
$ws.Cells.Item(18, 3).Value = 39.92297817715019
$ws.Cells.Item(18, 4).Value = 12.59090909090909
$ws.Cells.Item(18, 5).Value = 9.865928703026542
$ws.Cells.Item(18, 6).Value = -13.98476933769847
$ws.Cells.Item(18, 7).Value = -14546943.36
$ws.Cells.Item(18, 8).Value = 34.04343314010323
$ws.Cells.Item(18, 9).Value = 1040199.02
$ws.Cells.Item(18, 10).Value = 46.83244654002757
$ws.Cells.Item(18, 11).Value = 0.07150636351965559
$ws.Cells.Item(19, 2).Value = 164
$ws.Cells.Item(19, 3).Value = 21.05263157894737
$ws.Cells.Item(19, 4).Value = 14.94444444444444
$ws.Cells.Item(19, 5).Value = 9.869182409472742
$ws.Cells.Item(19, 6).Value = -27.8026173868846
$ws.Cells.Item(19, 7).Value = -15860891.66
$ws.Cells.Item(19, 8).Value = 37.11839603736733
$ws.Cells.Item(19, 9).Value = 570481.96
$ws.Cells.Item(19, 10).Value = 25.68457129843301
$ws.Cells.Item(19, 11).Value = 0.03596783662791881
$ws.Cells.Item(20, 2).Value = 232
$ws.Cells.Item(20, 3).Value = 29.78177150192555
$ws.Cells.Item(20, 4).Value = 13.58333333333333
$ws.Cells.Item(20, 5).Value = 10.13194677376464
$ws.Cells.Item(20, 6).Value = -15.30417706108132
$ws.Cells.Item(20, 7).Value = -6916346.34
$ws.Cells.Item(20, 8).Value = 16.18595524658644
$ws.Cells.Item(20, 9).Value = 451925.4
$ws.Cells.Item(20, 10).Value = 20.34684875552043
$ws.Cells.Item(20, 11).Value = 0.06534163816903363
$ws.Cells.Item(21, 2).Value = 72
$ws.Cells.Item(21, 3).Value = 9.242618741976894
$ws.Cells.Item(21, 4).Value = 18.33333333333333
$ws.Cells.Item(21, 5).Value = 9.832493791663456
$ws.Cells.Item(21, 6).Value = -34.10927008754508
$ws.Cells.Item(21, 7).Value = -5406360.24
$ws.Cells.Item(21, 8).Value = 12.65221557594299
$ws.Cells.Item(21, 9).Value = 158501.2
$ws.Cells.Item(21, 10).Value = 7.136133406018992
$ws.Cells.Item(21, 11).Value = 0.0293175432201684

# ---- XAUUSD block: refreshed header row 23 and data rows ----
$ws.Cells.Item(23, 2).Value = "count_of_occurrences"
$ws.Cells.Item(23, 3).Value = "percentage_of_occurrences"
$ws.Cells.Item(23, 4).Value = "typical_spread_in_points"
$ws.Cells.Item(23, 5).Value = "volume_weighted_avg_spread_in_USD"
$ws.Cells.Item(23, 6).Value = "PnL_per_lot"
$ws.Cells.Item(23, 7).Value = "total_profit"
$ws.Cells.Item(23, 8).Value = "pct_total_profit"
$ws.Cells.Item(23, 9).Value = "total_volume"
$ws.Cells.Item(23, 10).Value = "pct_total_volume"
$ws.Cells.Item(23, 11).Value = "pct_impact_on_PnL_exec_spread"
$ws.Cells.Item(25, 2).Value = 282
$ws.Cells.Item(25, 3).Value = 36.43410852713178
$ws.Cells.Item(25, 4).Value = 21.10204081632653
$ws.Cells.Item(25, 5).Value = 27.64771644344644
$ws.Cells.Item(25, 6).Value = -25.61993280568459
$ws.Cells.Item(25, 7).Value = -97662066.57
$ws.Cells.Item(25, 8).Value = 24.7580953969749
$ws.Cells.Item(25, 9).Value = 3811956.39
$ws.Cells.Item(25, 10).Value = 31.90527569151477
$ws.Cells.Item(25, 11).Value = 0.03903210861576181
$ws.Cells.Item(26, 2).Value = 132
$ws.Cells.Item(26, 3).Value = 17.05426356589147
$ws.Cells.Item(26, 4).Value = 19.5
$ws.Cells.Item(26, 5).Value = 31.27024985820477
$ws.Cells.Item(26, 6).Value = -67.71177022907793
$ws.Cells.Item(26, 7).Value = -145153531.57
$ws.Cells.Item(26, 8).Value = 36.79755208991035
$ws.Cells.Item(26, 9).Value = 2143697.19
$ws.Cells.Item(26, 10).Value = 17.94229598887817
$ws.Cells.Item(26, 11).Value = 0.0147684811166045
$ws.Cells.Item(27, 2).Value = 254
$ws.Cells.Item(27, 3).Value = 32.81653746770026
$ws.Cells.Item(27, 4).Value = 20.83333333333333
$ws.Cells.Item(27, 5).Value = 27.64915505976944
$ws.Cells.Item(27, 6).Value = -24.04970683129639
$ws.Cells.Item(27, 7).Value = -99321905.66
$ws.Cells.Item(27, 8).Value = 25.17887754891096
$ws.Cells.Item(27, 9).Value = 4129859.31
$ws.Cells.Item(27, 10).Value = 34.56605647388295
$ws.Cells.Item(27, 11).Value = 0.04158054844555024
$ws.Cells.Item(28, 2).Value = 106
$ws.Cells.Item(28, 3).Value = 13.69509043927649
$ws.Cells.Item(28, 4).Value = 21.92307692307692
$ws.Cells.Item(28, 5).Value = 30.11645615719254
$ws.Cells.Item(28, 6).Value = -28.09965755776398
$ws.Cells.Item(28, 7).Value = -52327680.23
$ws.Cells.Item(28, 8).Value = 13.26547496420377
$ws.Cells.Item(28, 9).Value = 1862217.72
$ws.Cells.Item(28, 10).Value = 15.58637184572409
$ws.Cells.Item(28, 11).Value = 0.03558762230266748

